$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values stay textually the same; shared string
# table reordering happens automatically since we write by text/value.

# Overwrite rows 2-29 with the new dataset (28 data rows).
$ws.Cells.Item(2, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(2, 2).Value = "402号直流"
$ws.Cells.Item(2, 3).Value = 45915.503680555557
$ws.Cells.Item(2, 4).Value = 45929.387083333335
$ws.Cells.Item(2, 5).Value = 333.20166666666046
$ws.Cells.Item(3, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(3, 2).Value = "101号直流"
$ws.Cells.Item(3, 3).Value = 45926.043692129628
$ws.Cells.Item(3, 4).Value = 45929.387083333335
$ws.Cells.Item(3, 5).Value = 80.241388888971414
$ws.Cells.Item(4, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(4, 2).Value = "603号直流"
$ws.Cells.Item(4, 3).Value = 45926.099652777775
$ws.Cells.Item(4, 4).Value = 45929.387083333335
$ws.Cells.Item(4, 5).Value = 78.898333333432674
$ws.Cells.Item(5, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(5, 2).Value = "602号直流"
$ws.Cells.Item(5, 3).Value = 45926.242071759261
$ws.Cells.Item(5, 4).Value = 45929.387083333335
$ws.Cells.Item(5, 5).Value = 75.480277777765878
$ws.Cells.Item(6, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(6, 2).Value = "008B号直流"
$ws.Cells.Item(6, 3).Value = 45926.525636574072
$ws.Cells.Item(6, 4).Value = 45929.387083333335
$ws.Cells.Item(6, 5).Value = 68.67472222232027
$ws.Cells.Item(7, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(7, 2).Value = "003B号直流"
$ws.Cells.Item(7, 3).Value = 45927.302870370368
$ws.Cells.Item(7, 4).Value = 45929.387083333335
$ws.Cells.Item(7, 5).Value = 50.02111111121485
$ws.Cells.Item(8, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(8, 2).Value = "701号直流"
$ws.Cells.Item(8, 3).Value = 45927.457337962966
$ws.Cells.Item(8, 4).Value = 45929.387083333335
$ws.Cells.Item(8, 5).Value = 46.313888888864312
$ws.Cells.Item(9, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(9, 2).Value = "705号直流"
$ws.Cells.Item(9, 3).Value = 45927.65315972222
$ws.Cells.Item(9, 4).Value = 45929.387083333335
$ws.Cells.Item(9, 5).Value = 41.61416666675359
$ws.Cells.Item(10, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10, 2).Value = "301号直流"
$ws.Cells.Item(10, 3).Value = 45928.017442129632
$ws.Cells.Item(10, 4).Value = 45929.387083333335
$ws.Cells.Item(10, 5).Value = 32.871388888859656
$ws.Cells.Item(11, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(11, 2).Value = "A01号直流"
$ws.Cells.Item(11, 3).Value = 45928.110659722224
$ws.Cells.Item(11, 4).Value = 45929.387083333335
$ws.Cells.Item(11, 5).Value = 30.634166666655801
$ws.Cells.Item(12, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(12, 2).Value = "001B号直流"
$ws.Cells.Item(12, 3).Value = 45928.127974537034
$ws.Cells.Item(12, 4).Value = 45929.387083333335
$ws.Cells.Item(12, 5).Value = 30.218611111224163
$ws.Cells.Item(13, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(13, 2).Value = "604号直流"
$ws.Cells.Item(13, 3).Value = 45928.226944444446
$ws.Cells.Item(13, 4).Value = 45929.387083333335
$ws.Cells.Item(13, 5).Value = 27.843333333323244
$ws.Cells.Item(14, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(14, 2).Value = "702号直流"
$ws.Cells.Item(14, 3).Value = 45928.532187500001
$ws.Cells.Item(14, 4).Value = 45929.387083333335
$ws.Cells.Item(14, 5).Value = 20.517500000016298
$ws.Cells.Item(15, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(15, 2).Value = "204号直流"
$ws.Cells.Item(15, 3).Value = 45928.540798611109
$ws.Cells.Item(15, 4).Value = 45929.387083333335
$ws.Cells.Item(15, 5).Value = 20.310833333409391
$ws.Cells.Item(16, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(16, 2).Value = "505号直流"
$ws.Cells.Item(16, 3).Value = 45928.549571759257
$ws.Cells.Item(16, 4).Value = 45929.387083333335
$ws.Cells.Item(16, 5).Value = 20.100277777877636
$ws.Cells.Item(17, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(17, 2).Value = "803号直流"
$ws.Cells.Item(17, 3).Value = 45928.553310185183
$ws.Cells.Item(17, 4).Value = 45929.387083333335
$ws.Cells.Item(17, 5).Value = 20.010555555636529
$ws.Cells.Item(18, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18, 2).Value = "404号直流"
$ws.Cells.Item(18, 3).Value = 45928.554293981484
$ws.Cells.Item(18, 4).Value = 45929.387083333335
$ws.Cells.Item(18, 5).Value = 19.98694444441935
$ws.Cells.Item(19, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(19, 2).Value = "210号直流"
$ws.Cells.Item(19, 3).Value = 45928.554675925923
$ws.Cells.Item(19, 4).Value = 45929.387083333335
$ws.Cells.Item(19, 5).Value = 19.977777777879965
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "502号直流"
$ws.Cells.Item(20, 3).Value = 45928.556516203702
$ws.Cells.Item(20, 4).Value = 45929.387083333335
$ws.Cells.Item(20, 5).Value = 19.933611111191567
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(21, 2).Value = "102号直流"
$ws.Cells.Item(21, 3).Value = 45928.562118055554
$ws.Cells.Item(21, 4).Value = 45929.387083333335
$ws.Cells.Item(21, 5).Value = 19.799166666751262
$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value = "901号直流"
$ws.Cells.Item(22, 3).Value = 45928.572928240741
$ws.Cells.Item(22, 4).Value = 45929.387083333335
$ws.Cells.Item(22, 5).Value = 19.539722222252749
$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23, 2).Value = "B02号直流"
$ws.Cells.Item(23, 3).Value = 45928.576990740738
$ws.Cells.Item(23, 4).Value = 45929.387083333335
$ws.Cells.Item(23, 5).Value = 19.442222222336568
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "306号直流"
$ws.Cells.Item(24, 3).Value = 45928.598923611113
$ws.Cells.Item(24, 4).Value = 45929.387083333335
$ws.Cells.Item(24, 5).Value = 18.915833333332557
$ws.Cells.Item(25, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(25, 2).Value = "107号直流"
$ws.Cells.Item(25, 3).Value = 45928.633726851855
$ws.Cells.Item(25, 4).Value = 45929.387083333335
$ws.Cells.Item(25, 5).Value = 18.080555555527098
$ws.Cells.Item(26, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(26, 2).Value = "111号直流"
$ws.Cells.Item(26, 3).Value = 45928.635104166664
$ws.Cells.Item(26, 4).Value = 45929.387083333335
$ws.Cells.Item(26, 5).Value = 18.047500000102445
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(27, 2).Value = "905号直流"
$ws.Cells.Item(27, 3).Value = 45928.713287037041
$ws.Cells.Item(27, 4).Value = 45929.387083333335
$ws.Cells.Item(27, 5).Value = 16.17111111106351
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(28, 2).Value = "B03号直流"
$ws.Cells.Item(28, 3).Value = 45928.737696759257
$ws.Cells.Item(28, 4).Value = 45929.387083333335
$ws.Cells.Item(28, 5).Value = 15.585277777863666
$ws.Cells.Item(29, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(29, 2).Value = "903号直流"
$ws.Cells.Item(29, 3).Value = 45928.766284722224
$ws.Cells.Item(29, 4).Value = 45929.387083333335
$ws.Cells.Item(29, 5).Value = 14.899166666669771

# Remove now-obsolete trailing rows 30-38 (old sheet had 38 rows, new has 29).
$ws.Rows("30:38").Delete()

# Move selection to match the authored state.
$ws.Range("H6").Select()
